$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2-481). Bump each value from 45189 (2023-09-20) to 45190
# (2023-09-21), matching the original spreadsheet's daily automatic
# update while leaving every other cell untouched.
$ws.Range("C2:C481").Value = 45190
